# ContosoLearn Market Research — translate body copy from Indonesian to
# English (matching the English-language source content) and collapse the
# multi-run paragraphs produced by the old generator into single runs.

$d = $word.ActiveDocument

function Set-ParagraphText($paraIndex, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    # Exclude the trailing paragraph mark so only the paragraph's own
    # content is replaced (this also collapses multiple runs into one).
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

Set-ParagraphText 1 "ContosoLearn Market Research"

Set-ParagraphText 2 "AdatumLearn: AdatumLearn is a top AI-powered learning platform that uses artificial intelligence to enrich eLearning with features that automate a variety of tasks. It is known for its content authoring capabilities and adaptive learning technology."

Set-ParagraphText 3 "AdventureLearn: AdventureLearn is another AI-powered learning platform that offers personalized learning experiences and data-driven recommendations."

Set-ParagraphText 4 "AlpineTraining: AlpineTraining is a mobile-first learning platform that focuses on microlearning."

Set-ParagraphText 5 "Bellows OnDemand: Bellows OnDemand is a comprehensive learning solution that offers content creation and social collaboration."

Set-ParagraphText 6 "FabrikamLearning: FabrikamLearning provides a suite of learning platforms that cater to different learning needs."

Set-ParagraphText 7 "FirstUp Cards: FirstUp Cards is a mobile learning app that is ideal for training on safety procedures, compliance, new product knowledge or any other type of training scenario."

Set-ParagraphText 8 "Munson'sLearn: Munson'sLearn is designed to enable businesses to train their employees, partners, and customers."

Set-ParagraphText 9 "LibertyLearn: LibertyLearn is a fast LMS for your mission-critical project."

Set-ParagraphText 10 "WoodgroveLMS: WoodgroveLMS is a functional and attractive learning management system built to provide a best-in-class training experience."

# The "a best" phrase is flagged by Word's grammar checker (article /
# adjective order) and ends up isolated in its own run, bracketed by the
# proofing-error range. Toggling a character property across just that
# span (scoped to the WoodgroveLMS paragraph) and back forces Word to
# split the run there instead of keeping the whole sentence as one run.
$grammarRange = $d.Paragraphs.Item(10).Range
$grammarRange.Find.Execute("a best", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$grammarRange.Bold = 1
$grammarRange.Bold = 0

Set-ParagraphText 11 "NorthwindWorlds: NorthwindWorlds is a powerful, easy-to-use, and reliable training solution for individuals and enterprises."

Set-ParagraphText 12 "ProsewareLearn: ProsewareLearn is an online education company that offers a variety of video training courses for software developers, IT administrators, and creative professionals through its website."

Set-ParagraphText 13 "RelecloudLearn: RelecloudLearn is an American online learning platform that offers massive open online courses (MOOC), specializations, and degrees in a variety of subjects."

Set-ParagraphText 14 "TreyAcademy: TreyAcademy is an online learning platform aimed at professional adults and students, developed in May 2010."

Set-ParagraphText 15 "These platforms have a significant market presence and are widely recognized for their AI-powered features, such as personalized learning experiences, data-driven recommendations, and automation of tasks. They are transforming the eLearning landscape by leveraging AI to deliver more engaging, rewarding, and personalized learning experiences. "
